$wb = $excel.ActiveWorkbook

# --- Rename the "Requested quantity" headers on the existing sheets ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Range("B1").Value = "Weekly_PO_Qty"

$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after "Monthly Trend" ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "PO Forecast"

# Header row
$ws3.Cells.Item(1, 1).Value = "ds"
$ws3.Cells.Item(1, 2).Value = "PO_Forecast"
$ws3.Cells.Item(1, 3).Value = "yhat_lower"
$ws3.Cells.Item(1, 4).Value = "yhat_upper"

$hdr = $ws3.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Data rows
$ws3.Cells.Item(2, 1).Value = 44934.99999999999
$ws3.Cells.Item(2, 2).Value = 84
$ws3.Cells.Item(2, 3).Value = 42.3554844641898
$ws3.Cells.Item(2, 4).Value = 128.0656007869159

$ws3.Cells.Item(3, 1).Value = 44969.99999999999
$ws3.Cells.Item(3, 2).Value = 61
$ws3.Cells.Item(3, 3).Value = 14.89897170040926
$ws3.Cells.Item(3, 4).Value = 102.4571800528219

$ws3.Cells.Item(4, 1).Value = 44976.99999999999
$ws3.Cells.Item(4, 2).Value = 56
$ws3.Cells.Item(4, 3).Value = 11.66457553834875
$ws3.Cells.Item(4, 4).Value = 99.09307366627218

$ws3.Cells.Item(5, 1).Value = 44983.99999999999
$ws3.Cells.Item(5, 2).Value = 52
$ws3.Cells.Item(5, 3).Value = 5.677109678984254
$ws3.Cells.Item(5, 4).Value = 95.0068313307895

$ws3.Cells.Item(6, 1).Value = 44990.99999999999
$ws3.Cells.Item(6, 2).Value = 47
$ws3.Cells.Item(6, 3).Value = 2.253350596585248
$ws3.Cells.Item(6, 4).Value = 92.62194988987351

$ws3.Cells.Item(7, 1).Value = 44997.99999999999
$ws3.Cells.Item(7, 2).Value = 43
$ws3.Cells.Item(7, 3).Value = -3.843812747666373
$ws3.Cells.Item(7, 4).Value = 86.18763981876747

$ws3.Cells.Item(8, 1).Value = 45004.99999999999
$ws3.Cells.Item(8, 2).Value = 38
$ws3.Cells.Item(8, 3).Value = -1.779864272246987
$ws3.Cells.Item(8, 4).Value = 83.37454048033788

$ws3.Cells.Item(9, 1).Value = 45011.99999999999
$ws3.Cells.Item(9, 2).Value = 33
$ws3.Cells.Item(9, 3).Value = -6.93084235636259
$ws3.Cells.Item(9, 4).Value = 78.35376932423519

$ws3.Cells.Item(10, 1).Value = 45018.99999999999
$ws3.Cells.Item(10, 2).Value = 29
$ws3.Cells.Item(10, 3).Value = -14.74791262553402
$ws3.Cells.Item(10, 4).Value = 71.56511685828839

$ws3.Cells.Item(11, 1).Value = 45025.99999999999
$ws3.Cells.Item(11, 2).Value = 24
$ws3.Cells.Item(11, 3).Value = -19.54540760568649
$ws3.Cells.Item(11, 4).Value = 70.22817445135561

$ws3.Cells.Item(12, 1).Value = 45032.99999999999
$ws3.Cells.Item(12, 2).Value = 19
$ws3.Cells.Item(12, 3).Value = -20.77057933109689
$ws3.Cells.Item(12, 4).Value = 65.42982344779267

$ws3.Cells.Item(13, 1).Value = 45039.99999999999
$ws3.Cells.Item(13, 2).Value = 15
$ws3.Cells.Item(13, 3).Value = -28.8999339051289
$ws3.Cells.Item(13, 4).Value = 58.15431461422837

$ws3.Cells.Item(14, 1).Value = 45046.99999999999
$ws3.Cells.Item(14, 2).Value = 10
$ws3.Cells.Item(14, 3).Value = -36.42489968573205
$ws3.Cells.Item(14, 4).Value = 56.52802532787778

# Apply the date number format (matching columns A on the other sheets) to column A data rows
$dateFmt = $ws1.Range("A2").NumberFormat
$ws3.Range("A2:A14").NumberFormat = $dateFmt

Write-Host "PO Forecast sheet created with" $ws3.UsedRange.Rows.Count "rows"
